$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$src = $ws.Range("F295:N295")
for ($r = 303; $r -le 311; $r++) {
    $dst = $ws.Range("F" + $r + ":N" + $r)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}
